# Updated remaining queries for C3DC
# - Rewrite the JOIN conditions in every SQL query cell so that the
#   generic ".id" columns are replaced with the fully-qualified
#   "<table>_id" columns (std.id -> std.study_id, prt.id -> prt.participant_id, ...).
# - Resize column C and drop its "best fit" auto-sizing.
# - Reset the sheet scroll position and move the active selection to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Query([string]$cellAddress) {
    $cell = $ws.Range($cellAddress)
    $text = $cell.Value()
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    $cell.Value = $text
}

# StatQuery (StudiesTab row) + TabQuery column for every tab.
Update-Query "C2"
Update-Query "B2"
Update-Query "B3"
Update-Query "B4"
Update-Query "B5"
Update-Query "B6"
Update-Query "B7"

# Column C: widen and switch off "best fit" sizing.
$ws.Columns.Item(3).ColumnWidth = 68.33

# Reset view: scroll back to the top-left and select B2.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B2").Select()
